$d = $word.ActiveDocument

# Locate the target paragraph: the liquidation announcement paragraph that
# still uses the short "sh." alias instead of "shareholders[0]." .
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*sh.street_number*") {
        $target = $para
        break
    }
}

$pStart = $target.Range.Start

# Turn on revision tracking so that replacing the four "sh." occurrences
# creates genuine new runs (rather than being silently coalesced back into
# the single surrounding run). We then accept all the revisions, which
# keeps the run split produced by the edit while dropping the ins/del
# wrapper markup, exactly mirroring how Word materializes multiple
# small in-place replacements as separate <w:r> runs.
$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true

# Replace occurrences from last to first so earlier offsets stay valid.
$offsets = @(372, 354, 333, 310)
foreach ($off in $offsets) {
    $r = $d.Range($pStart + $off, $pStart + $off + 3)
    $r.Text = "shareholders[0]."
}

$d.TrackRevisions = $wasTracking

# Accept revisions one at a time (instead of Revisions.AcceptAll(), which
# as a side effect strips the w:rsidRPr attribute off every run in the
# whole document, not just the ones touched by this edit).
while ($d.Revisions.Count -gt 0) {
    $d.Revisions.Item(1).Accept()
}

Write-Output $target.Range.Text
